# Update view-count ("想去人数", column F) and lowest-price ("最低票价",
# column G) figures across the four sheets, matching the freshly generated
# gh-pages data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 511
$ws1.Range("F3").Value  = 1593
$ws1.Range("F4").Value  = 841
$ws1.Range("F5").Value  = 239
$ws1.Range("F6").Value  = 69
$ws1.Range("F7").Value  = 1147
$ws1.Range("F9").Value  = 795
$ws1.Range("F10").Value = 1446
$ws1.Range("F12").Value = 1035
$ws1.Range("F18").Value = 26
$ws1.Range("F19").Value = 26
$ws1.Range("F23").Value = 553
$ws1.Range("F24").Value = 569
$ws1.Range("F25").Value = 759
$ws1.Range("F26").Value = 246
$ws1.Range("F27").Value = 179

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
# Row 9's event is no longer sold out, so it now reports a numeric
# lowest price (580) instead of the "已售罄" ("sold out") text.
$ws2.Range("F9").Value  = 587
$ws2.Range("G9").Value  = 580
$ws2.Range("F11").Value = 13

# --- Sheet "本地生活" (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 237

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 511
$ws4.Range("F3").Value  = 237
$ws4.Range("F4").Value  = 1593
$ws4.Range("F6").Value  = 841
$ws4.Range("F7").Value  = 239
$ws4.Range("F9").Value  = 69
$ws4.Range("F10").Value = 1147
$ws4.Range("F12").Value = 795
$ws4.Range("F13").Value = 1446
$ws4.Range("F15").Value = 1035
$ws4.Range("F21").Value = 26
$ws4.Range("F22").Value = 26
$ws4.Range("F31").Value = 553
$ws4.Range("F32").Value = 569
$ws4.Range("F33").Value = 759
$ws4.Range("F34").Value = 246
$ws4.Range("F36").Value = 179
# Row 37 mirrors sheet "演出" row 9 above.
$ws4.Range("F37").Value = 587
$ws4.Range("G37").Value = 580
$ws4.Range("F40").Value = 13
